# bm25 run and reranking run (near finished)
# Replace the "x" marker (shared string) cells in rows 19-26 with plain
# numeric rank/run values, and move the active selection/viewport.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19 -> 0
foreach ($addr in @("S19","U19","AD19","AE19","AF19","AG19","AL19","AM19")) {
    $ws.Range($addr).Value = 0
}

# Row 20 -> 1
foreach ($addr in @("AK20","AN20","AO20","AP20")) {
    $ws.Range($addr).Value = 1
}

# Row 21 -> 2
foreach ($addr in @("R21")) {
    $ws.Range($addr).Value = 2
}

# Row 22 -> 3
foreach ($addr in @("T22","V22","X22","AS22","AT22","AU22")) {
    $ws.Range($addr).Value = 3
}

# Row 23 -> 4
foreach ($addr in @("AB23")) {
    $ws.Range($addr).Value = 4
}

# Row 24 -> 5
foreach ($addr in @("AA24","AC24","AI24","AJ24")) {
    $ws.Range($addr).Value = 5
}

# Row 25 -> 6
foreach ($addr in @("Z25")) {
    $ws.Range($addr).Value = 6
}

# Row 26 -> 7
foreach ($addr in @("W26","Y26","AH26","AQ26","AR26")) {
    $ws.Range($addr).Value = 7
}

# Move viewport / active cell to match the saved view state.
$win = $excel.ActiveWindow
$win.ScrollColumn = 27
$win.ScrollRow = 1
$ws.Range("AS22").Select()
